$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USD Amount value in T2
$ws.Range("T2").Value = 216399

# Update the selected/active cell on the sheet (also clears the old
# topLeftCell scroll position, since the sheet view recenters on selection)
$ws.Range("H22").Select()
